$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header text updates (Volume/Number and report-week date range) ---
# Both shared strings are built from multiple identically-formatted runs;
# Value assignment rewrites the whole string but renders identically.
$ws.Range("A8").Value = "Volume 31   Number  42"
$ws.Range("C9").Value = "Report Covering the Week  10/14/2024  Through  10/20/2024"

# --- Weekly crime-stat grid updates (rows 14-30) ---
# Some cells flip between the text "N/A" placeholder style (used when a 2023
# count was 0, to avoid dividing by zero) and a real numeric value. Plain
# Value assignment cannot change a text-style cell to a genuine number style
# (or vice-versa) without altering the cell style id, so for those specific
# cells we first clone the correct format via Range.Copy from a same-column
# donor cell of the desired type, then (for text targets) the copy already
# carries the right literal value, or (for numeric targets) we overwrite the
# value afterwards.

# --- Row 14 ---
$ws.Range("L14").Value = -62.5

# --- Row 15 ---
$ws.Range("F15").Value = 2
$ws.Range("G15").Value = 1
$ws.Range("H15").Value = 100
$ws.Range("L15").Value = 25
$ws.Range("M15").Value = 53.846153846153

# --- Row 16 ---
$ws.Range("F16").Value = 10
$ws.Range("H16").Value = 25
$ws.Range("I16").Value = 110
$ws.Range("J16").Value = 85
$ws.Range("K16").Value = 29.411764705882
$ws.Range("L16").Value = 46.666666666666
$ws.Range("M16").Value = 11.111111111111
$ws.Range("N16").Value = -79.92700729927

# --- Row 17 ---
$ws.Range("C17").Value = 11
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = 120
$ws.Range("F17").Value = 28
$ws.Range("G17").Value = 23
$ws.Range("H17").Value = 21.739130434782
$ws.Range("I17").Value = 287
$ws.Range("J17").Value = 218
$ws.Range("K17").Value = 31.651376146789
$ws.Range("L17").Value = 51.052631578947
$ws.Range("M17").Value = 124.21875
$ws.Range("N17").Value = -30.843373493975

# --- Row 18 ---
$ws.Range("C16").Copy($ws.Range("C18"))
$ws.Range("C18").Value = 3
$ws.Range("E18").Value = 200
$ws.Range("F18").Value = 7
$ws.Range("G18").Value = 4
$ws.Range("H18").Value = 75
$ws.Range("I18").Value = 59
$ws.Range("J18").Value = 57
$ws.Range("K18").Value = 3.508771929824
$ws.Range("L18").Value = -3.27868852459
$ws.Range("M18").Value = -45.871559633027
$ws.Range("N18").Value = -90.895061728395

# --- Row 19 ---
$ws.Range("C19").Value = 2
$ws.Range("D19").Value = 2
$ws.Range("G19").Value = 15
$ws.Range("H19").Value = 6.666666666666
$ws.Range("I19").Value = 135
$ws.Range("J19").Value = 156
$ws.Range("K19").Value = -13.461538461538
$ws.Range("L19").Value = -5.594405594405
$ws.Range("M19").Value = 70.886075949367
$ws.Range("N19").Value = -48.669201520912

# --- Row 20 ---
$ws.Range("C20").Value = 2
$ws.Range("C14").Copy($ws.Range("D20"))
$ws.Range("E14").Copy($ws.Range("E20"))
$ws.Range("G20").Value = 4
$ws.Range("H20").Value = 100
$ws.Range("I20").Value = 64
$ws.Range("K20").Value = 20.754716981132
$ws.Range("L20").Value = 25.490196078431
$ws.Range("M20").Value = -5.882352941176
$ws.Range("N20").Value = -82.933333333333

# --- Row 21 ---
$ws.Range("C21").Value = 20
$ws.Range("D21").Value = 10
$ws.Range("E21").Value = 100
$ws.Range("F21").Value = 72
$ws.Range("G21").Value = 55
$ws.Range("H21").Value = 30.90909090909
$ws.Range("I21").Value = 678
$ws.Range("J21").Value = 582
$ws.Range("K21").Value = 16.494845360824
$ws.Range("L21").Value = 24.632352941176
$ws.Range("M21").Value = 35.329341317365
$ws.Range("N21").Value = -70.5089169204

# --- Row 22 ---
$ws.Range("L22").Value = 25

# --- Row 23 ---
$ws.Range("C16").Copy($ws.Range("C23"))
$ws.Range("C23").Value = 2
$ws.Range("D16").Copy($ws.Range("D23"))
$ws.Range("D23").Value = 1
$ws.Range("E16").Copy($ws.Range("E23"))
$ws.Range("E23").Value = 100
$ws.Range("F23").Value = 7
$ws.Range("G23").Value = 3
$ws.Range("H23").Value = 133.333333333333
$ws.Range("I23").Value = 70
$ws.Range("J23").Value = 67
$ws.Range("K23").Value = 4.477611940298
$ws.Range("L23").Value = -1.408450704225
$ws.Range("M23").Value = 105.882352941176

# --- Row 24 ---
$ws.Range("C24").Value = 11
$ws.Range("D24").Value = 9
$ws.Range("E24").Value = 22.222222222222
$ws.Range("F24").Value = 60
$ws.Range("G24").Value = 42
$ws.Range("H24").Value = 42.857142857142
$ws.Range("I24").Value = 537
$ws.Range("J24").Value = 466
$ws.Range("K24").Value = 15.236051502145
$ws.Range("L24").Value = 12.81512605042
$ws.Range("M24").Value = 80.80808080808

# --- Row 25 ---
$ws.Range("C14").Copy($ws.Range("C25"))
$ws.Range("E25").Value = -100
$ws.Range("F25").Value = 2
$ws.Range("G25").Value = 9
$ws.Range("H25").Value = -77.777777777777
$ws.Range("J25").Value = 81
$ws.Range("K25").Value = -53.086419753086
$ws.Range("L25").Value = -47.222222222222

# --- Row 26 ---
$ws.Range("C26").Value = 11
$ws.Range("D26").Value = 10
$ws.Range("E26").Value = 10
$ws.Range("F26").Value = 42
$ws.Range("G26").Value = 30
$ws.Range("H26").Value = 40
$ws.Range("I26").Value = 354
$ws.Range("J26").Value = 346
$ws.Range("K26").Value = 2.312138728323
$ws.Range("L26").Value = 9.937888198757
$ws.Range("M26").Value = 5.357142857142

# --- Row 27 ---
$ws.Range("F27").Value = 2
$ws.Range("G27").Value = 1
$ws.Range("H27").Value = 100
$ws.Range("L27").Value = 6.896551724137

# --- Row 28 ---
$ws.Range("C16").Copy($ws.Range("C28"))
$ws.Range("C28").Value = 2
$ws.Range("F28").Value = 5
$ws.Range("G28").Value = 1
$ws.Range("H28").Value = 400
$ws.Range("I28").Value = 37
$ws.Range("K28").Value = 8.823529411764
$ws.Range("L28").Value = -15.90909090909

# --- Row 29 ---
$ws.Range("F29").Value = 1
$ws.Range("H29").Value = -66.666666666666
$ws.Range("L29").Value = -40

# --- Row 30 ---
$ws.Range("F30").Value = 1
$ws.Range("H30").Value = -66.666666666666
$ws.Range("L30").Value = -31.578947368421
